$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 89) below the existing data (which ends at row 88).
# Column A holds a date written as plain text (e.g. "2025/10/10"), matching the
# existing rows which store dates as text rather than native date values.
# Copying the cell above (same date text, default/no style) avoids Excel's
# automatic text->date conversion and keeps the same formatting (none).
$ws.Range("A88").Copy($ws.Range("A89"))

$ws.Range("B89").Value = "金"
$ws.Range("C89").Value = 16
$ws.Range("D89").Value = 32
